# "南宁-漫展信息.xlsx" update - a new event "南宁·恋与深空only" was added,
# and the "想去人数" (interest counts) for two existing events were refreshed.
#
# This touches two worksheets:
#   - "展览"   (sheet1): new event appended as row 5 (A1:I4 -> A1:I5)
#   - "全部类型" (sheet4): new event inserted as row 5, pushing the existing
#                         last row ("...浪漫古典...") down to row 6
#                         (A1:I5 -> A1:I6)

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$text) {
    # Force the cell to stay a text value even when the content looks like
    # a date (e.g. "2024-06-09"), then drop back to the default/no style
    # so no stray formatting is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

function Set-NewEventRow($ws, [int]$row, [int]$seq, [string]$copyStyleFromRow) {
    $ws.Cells.Item($copyStyleFromRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($row, 1).Value = $seq

    Set-TextValue $ws.Cells.Item($row, 2) "2024-06-09"
    Set-TextValue $ws.Cells.Item($row, 3) "南宁·恋与深空only"
    Set-TextValue $ws.Cells.Item($row, 4) "新阳路227号南宁第三人民医院旁新秀佳园对面 卡尔顿东方银龙酒店"
    Set-TextValue $ws.Cells.Item($row, 5) "2024.06.09 10:00-06.09 17:00"
    $ws.Cells.Item($row, 6).Value = 2
    $ws.Cells.Item($row, 7).Value = 50
    Set-TextValue $ws.Cells.Item($row, 8) "https://show.bilibili.com/platform/detail.html?id=84444"
    Set-TextValue $ws.Cells.Item($row, 9) "//i2.hdslb.com/bfs/openplatform/202404/6ZVHU1F91713340880421.jpeg"
}

# ---- Sheet "展览": refresh counts and append new row 5 ----
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value = 5506
$ws1.Cells.Item(4, 6).Value = 940

Set-NewEventRow $ws1 5 4 4

# ---- Sheet "全部类型": refresh counts and insert the new row before the last one ----
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value = 5506
$ws4.Cells.Item(4, 6).Value = 940

$ws4.Rows.Item(5).Insert()

Set-NewEventRow $ws4 5 4 4

# Renumber the event that got pushed down from row 5 to row 6
$ws4.Cells.Item(6, 1).Value = 5
